$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 597.3333
$ws.Range("I18").Value = 795
$ws.Range("J18").Value = 202
$ws.Range("K18").Value = 795
$ws.Range("L18").Value = 202
$ws.Range("M18").Value = -511
$ws.Range("N18").Value = -770
$ws.Range("H40").Value = 13876.823
$ws.Range("I40").Value = 2981.4
$ws.Range("J40").Value = 18416.584
$ws.Range("K40").Value = 2981.4
$ws.Range("L40").Value = 18416.584
$ws.Range("M40").Value = -2806.4
$ws.Range("N40").Value = -18766.584
$ws.Range("H64").Value = 7536
$ws.Range("I64").Value = 6816.1816
$ws.Range("J64").Value = 8667.143
$ws.Range("K64").Value = 6816.1816
$ws.Range("L64").Value = 8667.143
$ws.Range("M64").Value = -6568.1816
$ws.Range("N64").Value = -9163.143
$ws.Range("H67").Value = 7536
$ws.Range("I67").Value = 6816.1816
$ws.Range("J67").Value = 8667.143
$ws.Range("K67").Value = 6816.1816
$ws.Range("L67").Value = 8667.143
$ws.Range("M67").Value = -5958.1816
$ws.Range("N67").Value = -10383.143
$ws.Range("H115").Value = 2260.3333
$ws.Range("I115").Value = 2260.3333
$ws.Range("K115").Value = 6780.999899999999
$ws.Range("M115").Value = -5213.999899999999
$ws.Range("H123").Value = 89999
$ws.Range("J123").Value = 89999
$ws.Range("L123").Value = 89999
$ws.Range("N123").Value = -99799
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3491.7778
$ws.Range("I2").Value = 2715.3
$ws.Range("K2").Value = 2715.3
$ws.Range("M2").Value = -2602.3
$ws.Range("H45").Value = 3317.76
$ws.Range("I45").Value = 2246.8125
$ws.Range("J45").Value = 5221.6665
$ws.Range("K45").Value = 2246.8125
$ws.Range("L45").Value = 5221.6665
$ws.Range("M45").Value = -1869.8125
$ws.Range("N45").Value = -5975.6665
$ws.Range("H63").Value = 6606.381
$ws.Range("I63").Value = 5775.5264
$ws.Range("J63").Value = 14499.5
$ws.Range("K63").Value = 5775.5264
$ws.Range("L63").Value = 14499.5
$ws.Range("M63").Value = -5089.5264
$ws.Range("N63").Value = -15871.5
$ws.Range("H66").Value = 6606.381
$ws.Range("I66").Value = 5775.5264
$ws.Range("J66").Value = 14499.5
$ws.Range("K66").Value = 28877.632
$ws.Range("L66").Value = 72497.5
$ws.Range("M66").Value = -25445.632
$ws.Range("N66").Value = -79361.5
$ws.Range("H116").Value = 3491.7778
$ws.Range("I116").Value = 2715.3
$ws.Range("K116").Value = 2715.3
$ws.Range("M116").Value = -421.3000000000002
$ws.Range("H122").Value = 4067.4
$ws.Range("I122").Value = 4067.4
$ws.Range("K122").Value = 12202.2
$ws.Range("M122").Value = -9752.200000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3491.7778
$ws.Range("I3").Value = 2715.3
$ws.Range("K3").Value = 2715.3
$ws.Range("M3").Value = -2601.3
$ws.Range("H86").Value = 3648.8572
$ws.Range("I86").Value = 2693
$ws.Range("J86").Value = 4365.75
$ws.Range("K86").Value = 2693
$ws.Range("L86").Value = 4365.75
$ws.Range("M86").Value = -1570
$ws.Range("N86").Value = -6611.75
$ws.Range("H89").Value = 3648.8572
$ws.Range("I89").Value = 2693
$ws.Range("J89").Value = 4365.75
$ws.Range("K89").Value = 13465
$ws.Range("L89").Value = 21828.75
$ws.Range("M89").Value = -7849
$ws.Range("N89").Value = -33060.75
$ws.Range("H107").Value = 2536.6316
$ws.Range("I107").Value = 2536.6316
$ws.Range("K107").Value = 2536.6316
$ws.Range("M107").Value = -616.6316000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6218.207
$ws.Range("I31").Value = 7103.4585
$ws.Range("J31").Value = 1969
$ws.Range("K31").Value = 7103.4585
$ws.Range("L31").Value = 1969
$ws.Range("M31").Value = -6808.4585
$ws.Range("N31").Value = -2559
$ws.Range("H34").Value = 6218.207
$ws.Range("I34").Value = 7103.4585
$ws.Range("J34").Value = 1969
$ws.Range("K34").Value = 7103.4585
$ws.Range("L34").Value = 1969
$ws.Range("M34").Value = -6901.4585
$ws.Range("N34").Value = -2373
$ws.Range("H58").Value = 3514.2307
$ws.Range("I58").Value = 3898.6365
$ws.Range("J58").Value = 1400
$ws.Range("K58").Value = 3898.6365
$ws.Range("L58").Value = 1400
$ws.Range("M58").Value = -3695.6365
$ws.Range("N58").Value = -1806
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H136").Value = 3514.2307
$ws.Range("I136").Value = 3898.6365
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 11695.9095
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -9145.9095
$ws.Range("N136").Value = -9300
$ws.Range("H141").Value = 81298.25
$ws.Range("J141").Value = 94633
$ws.Range("L141").Value = 94633
$ws.Range("N141").Value = -104993
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1444.8889
$ws.Range("I5").Value = 1625.2858
$ws.Range("J5").Value = 813.5
$ws.Range("K5").Value = 4875.857400000001
$ws.Range("L5").Value = 2440.5
$ws.Range("M5").Value = -4763.857400000001
$ws.Range("N5").Value = -2664.5
$ws.Range("H37").Value = 125059910
$ws.Range("J37").Value = 125059910
$ws.Range("L37").Value = 375179730
$ws.Range("N37").Value = -375179954
$ws.Range("H68").Value = 816.3333
$ws.Range("I68").Value = 699
$ws.Range("K68").Value = 2097
$ws.Range("M68").Value = -1286
$ws.Range("H71").Value = 816.3333
$ws.Range("I71").Value = 699
$ws.Range("K71").Value = 6291
$ws.Range("M71").Value = -2235
$ws.Range("H107").Value = 914.5714
$ws.Range("I107").Value = 519.6
$ws.Range("K107").Value = 1558.8
$ws.Range("M107").Value = 361.1999999999998
$ws.Range("H121").Value = 840.8889
$ws.Range("J121").Value = 1599.75
$ws.Range("L121").Value = 4799.25
$ws.Range("N121").Value = -7419.25
$ws.Range("H135").Value = 1444.8889
$ws.Range("I135").Value = 1625.2858
$ws.Range("J135").Value = 813.5
$ws.Range("K135").Value = 14627.5722
$ws.Range("L135").Value = 7321.5
$ws.Range("M135").Value = -12092.5722
$ws.Range("N135").Value = -12391.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 30342.666
$ws.Range("I47").Value = 25028
$ws.Range("K47").Value = 25028
$ws.Range("M47").Value = -24460
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H99").Value = 16042.375
$ws.Range("I99").Value = 8334.286
$ws.Range("J99").Value = 69999
$ws.Range("K99").Value = 8334.286
$ws.Range("L99").Value = 69999
$ws.Range("M99").Value = -6088.286
$ws.Range("N99").Value = -74491
$ws.Range("H119").Value = 80000
$ws.Range("J119").Value = 80000
$ws.Range("L119").Value = 80000
$ws.Range("N119").Value = -89676
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31760.3
$ws.Range("I7").Value = 29225.375
$ws.Range("K7").Value = 29225.375
$ws.Range("M7").Value = -29113.375
$ws.Range("H40").Value = 11612.435
$ws.Range("I40").Value = 8472.643
$ws.Range("K40").Value = 8472.643
$ws.Range("M40").Value = -8336.643
$ws.Range("H46").Value = 5937
$ws.Range("I46").Value = 6125
$ws.Range("K46").Value = 6125
$ws.Range("M46").Value = -5937
$ws.Range("H61").Value = 3247
$ws.Range("I61").Value = 3247
$ws.Range("K61").Value = 3247
$ws.Range("M61").Value = -3045
$ws.Range("H113").Value = 3247
$ws.Range("I113").Value = 3247
$ws.Range("K113").Value = 3247
$ws.Range("M113").Value = -1077
$ws.Range("H122").Value = 6416.222
$ws.Range("I122").Value = 2791.1667
$ws.Range("K122").Value = 8373.500100000001
$ws.Range("M122").Value = -5923.500100000001
$ws.Range("H126").Value = 31760.3
$ws.Range("I126").Value = 29225.375
$ws.Range("K126").Value = 87676.125
$ws.Range("M126").Value = -85206.125
$ws.Range("H132").Value = 4287.7896
$ws.Range("I132").Value = 4483.5
$ws.Range("K132").Value = 13450.5
$ws.Range("M132").Value = -10920.5
$ws.Range("H140").Value = 77729.38
$ws.Range("J140").Value = 72961
$ws.Range("L140").Value = 72961
$ws.Range("N140").Value = -83321
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 33660.266
$ws.Range("J81").Value = 73268.336
$ws.Range("L81").Value = 146536.672
$ws.Range("N81").Value = -148658.672
$ws.Range("H84").Value = 33660.266
$ws.Range("J84").Value = 73268.336
$ws.Range("L84").Value = 732683.36
$ws.Range("N84").Value = -743291.36
$ws.Range("H107").Value = 2008.3103
$ws.Range("I107").Value = 2465.75
$ws.Range("J107").Value = 1685.4117
$ws.Range("K107").Value = 7397.25
$ws.Range("L107").Value = 5056.2351
$ws.Range("M107").Value = -5477.25
$ws.Range("N107").Value = -8896.2351
$ws.Range("H113").Value = 685.5263
$ws.Range("I113").Value = 685.5263
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2056.5789
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 113.4211
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 3959.8
$ws.Range("I126").Value = 3959.8
$ws.Range("K126").Value = 11879.4
$ws.Range("M126").Value = -9409.400000000001
$ws.Range("H132").Value = 3120.353
$ws.Range("I132").Value = 3174.875
$ws.Range("K132").Value = 9524.625
$ws.Range("M132").Value = -6994.625
